$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grade sheet: fill in previously-blank attendance/grade cells with 5
$ws.Range("C10:F10").Value = 5
$ws.Range("C19").Value = 5
$ws.Range("C20").Value = 5
$ws.Range("C26").Value = 5
$ws.Range("C27:F27").Value = 5
$ws.Range("C30:F30").Value = 5

# Move the active selection to F10, matching the author's final cursor position
$ws.Range("F10").Select()
